# Update doctor records on the "Doctors" sheet (rows 2-6): replace the
# Name (A), Experience (C) and Practise Location (D) entries with a new
# batch of doctors/locations, per the commit "Removed JS for some, so
# that it runs for eclipse".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Doctors")

# Row 2
$ws.Range("A2").Value = "Dr. Shantha Rama Rao"
$ws.Range("C2").Value = "51 years experience overall"
$ws.Range("D2").Value = "Kodambakkam,Chennai"

# Row 3
$ws.Range("A3").Value = "Dr. Bobby M"
$ws.Range("C3").Value = "22 years experience overall"
$ws.Range("D3").Value = "Thoraipakkam,Chennai"

# Row 4
$ws.Range("A4").Value = "Dr. Sathya Balasubramanyam"
$ws.Range("C4").Value = "28 years experience overall"
$ws.Range("D4").Value = "T Nagar,Chennai"

# Row 5
$ws.Range("A5").Value = "Dr. V. Bharathi"
$ws.Range("C5").Value = "17 years experience overall"
$ws.Range("D5").Value = "Mandaveli,Chennai"

# Row 6
$ws.Range("A6").Value = "Dr. Karthiga Devi"
$ws.Range("C6").Value = "19 years experience overall"
$ws.Range("D6").Value = "Karapakkam,Chennai"
